$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new data for row 4 (Rotate_60_high), which previously had no
# measurement data (only the label in A4/F4).
$ws.Range("B4").Value = "3.45（4）"
$ws.Range("C4").Value = "742（35）"

# Add explanatory note in D10 about why chamfer conflicts with kernel.
$ws.Range("D10").Value = "为何chamfer于kernel犯冲：事件数量！Kernel方法受事件数量影响严重，要保证输入的事件数量尽量一致"

$ws.Range("G4").Value = "3.7022（5）"
$ws.Range("H4").Value = "1965（49s)"

# Update the selection / view state to match the post-edit workbook.
$ws.Range("F11").Select()
